# "add arrow figures for slideshow"
#
# 1) Bump the cached datetimeFigureOut placeholder text (slide master + every
#    slide layout) from 6/12/18 to 6/13/18.
# 2) Re-layout the single content slide: widen/centre/re-wrap the title
#    textbox so it spans the full slide width, and nudge every other shape
#    (the figure picture, labels, rectangles and connector arrows) left by
#    the same amount to make room for new arrow figures.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Date placeholder text: slide master + all custom (slide) layouts.
# ---------------------------------------------------------------------
function Set-DatePlaceholderText {
    param($shapes, $text)
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $isDate = $false
        try { if ($sh.PlaceholderFormat.Type -eq 16) { $isDate = $true } } catch {}
        if ($isDate) {
            $sh.TextFrame.TextRange.Text = $text
        }
    }
}

$master = $p.SlideMaster
Set-DatePlaceholderText $master.Shapes "6/13/18"

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Set-DatePlaceholderText $layout.Shapes "6/13/18"
}

# ---------------------------------------------------------------------
# 2) Slide 1 shape re-layout.
# ---------------------------------------------------------------------
$s = $p.Slides.Item(1)

function Get-ShapeById {
    param($slide, $id)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Id -eq $id) { return $sh }
    }
    return $null
}

# Title textbox ("Mass spectrometrists ...") - id 19: now spans the full
# slide width, centred, word-wrapped instead of a single non-wrapping line.
$title = Get-ShapeById $s 19
$title.Left = -0.00007874015864218
$title.Width = 850.37506103515625
$title.TextFrame.WordWrap = -1
$title.TextFrame.TextRange.ParagraphFormat.Alignment = 2

# All remaining shapes only shift horizontally (Left changes, nothing else).
$leftShifts = @{
    7  = 294.98553466796875    # Picture 6 (figure image)
    2  = 24.32685089111328125  # TextBox 1 ("Search against DB ...")
    5  = 28.72448921203613281  # TextBox 4 ("Remove irrelevant PSMs ...")
    6  = 29.69874191284179688  # TextBox 5 ("FDR calculation on subset ...")
    10 = 145.1656036376953125  # Straight Arrow Connector 9
    9  = 289.24481201171875    # Rectangle 8
    12 = 254.5230712890625     # Straight Connector 11
    14 = 255.497406005859375   # Straight Connector 13
    17 = 145.1656036376953125  # Straight Arrow Connector 16
    18 = 28.7244110107421875   # Rectangle 17
    21 = 28.7244110107421875   # Rectangle 20
    22 = 29.69874191284179688  # Rectangle 21
}

foreach ($id in $leftShifts.Keys) {
    $sh = Get-ShapeById $s $id
    $sh.Left = $leftShifts[$id]
}
